# Order of operations matters for shared-string interning order:
# ThrowKnife, then Mp, then Evade (matches the target sharedStrings.xml order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: ThrowKnife (B column string interned first -> index 6)
$ws.Range("B5").Value = "ThrowKnife"

# New "Mp" column header (F2) -> interned second -> index 7
$ws.Range("F2").Value = "Mp"

# Mp values for existing rows
$ws.Range("F3").Value = 15
$ws.Range("F4").Value = 30

# Finish row 5 values
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 10

# New row 6: Evade -> interned third -> index 8
$ws.Range("B6").Value = "Evade"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 10

# Update selection to match target state
$ws.Range("H4").Select()
